$wb = $excel.ActiveWorkbook

# --- Update selection on the "Company" sheet (was tabSelected, selection A3) ---
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Activate()
$wsCompany.Range("K24").Select()

# --- Add the new "Activity" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsActivity = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsActivity.Name = "Activity"

# Header row (bold + centered)
$wsActivity.Range("A1").Value = "Type"
$wsActivity.Range("B1").Value = "Subject"
$wsActivity.Range("C1").Value = "IndustryGroup"
$wsActivity.Range("D1").Value = "ProductType"
$wsActivity.Range("E1").Value = "Description"
$wsActivity.Range("F1").Value = "MeetingNotes"
$wsActivity.Range("G1").Value = "ExtAttendee"

$wsActivity.Range("A1:G1").Font.Bold = $true
$wsActivity.Range("A1:G1").HorizontalAlignment = -4108

# Data row
$wsActivity.Range("A2").Value = "Meeting"
$wsActivity.Range("B2").Value = "Test Conf 01"
$wsActivity.Range("C2").Value = "BUS - Business Services"
$wsActivity.Range("D2").Value = "Activist Advisory"
$wsActivity.Range("E2").Value = "Test Conf 01"
$wsActivity.Range("F2").Value = "Test Conf 01"
$wsActivity.Range("G2").Value = "Test External"

# Column widths (best-fit to content)
$wsActivity.Columns.Item(2).AutoFit()
$wsActivity.Columns.Item(3).AutoFit()
$wsActivity.Columns.Item(4).AutoFit()
$wsActivity.Columns.Item(5).AutoFit()
$wsActivity.Columns.Item(6).AutoFit()
$wsActivity.Columns.Item(7).AutoFit()

# Final selection/active state on the new sheet
$wsActivity.Range("H10").Select()
